$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 348.18182
$ws.Range("I2").Value = 305.3
$ws.Range("K2").Value = 305.3
$ws.Range("M2").Value = -192.3
$ws.Range("H17").Value = 995.3674
$ws.Range("J17").Value = 1069.8182
$ws.Range("L17").Value = 3209.4546
$ws.Range("N17").Value = -3545.4546
$ws.Range("H32").Value = 3842.4285
$ws.Range("H62").Value = 9548.615
$ws.Range("J62").Value = 10389.777
$ws.Range("L62").Value = 10389.777
$ws.Range("N62").Value = -11637.777
$ws.Range("H65").Value = 9548.615
$ws.Range("J65").Value = 10389.777
$ws.Range("L65").Value = 51948.885
$ws.Range("N65").Value = -58188.885
$ws.Range("H74").Value = 3175.889
$ws.Range("I74").Value = 3322.875
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 3322.875
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -2386.875
$ws.Range("N74").Value = -3872
$ws.Range("H77").Value = 3175.889
$ws.Range("I77").Value = 3322.875
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 16614.375
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -11934.375
$ws.Range("N77").Value = -19360
$ws.Range("H132").Value = 2828.92
$ws.Range("I132").Value = 2209.8667
$ws.Range("J132").Value = 3757.5
$ws.Range("K132").Value = 6629.6001
$ws.Range("L132").Value = 11272.5
$ws.Range("M132").Value = -4099.6001
$ws.Range("N132").Value = -16332.5
$ws.Range("H135").Value = 1301.5714
$ws.Range("I135").Value = 1375.1538
$ws.Range("J135").Value = 345
$ws.Range("K135").Value = 12376.3842
$ws.Range("L135").Value = 3105
$ws.Range("M135").Value = -9841.3842
$ws.Range("N135").Value = -8175
$ws.Range("H137").Value = 2133.3845
$ws.Range("I137").Value = 1438.7778
$ws.Range("K137").Value = 4316.3334
$ws.Range("M137").Value = -1766.3334
$ws.Range("H138").Value = 4603.28
$ws.Range("I138").Value = 925.7059
$ws.Range("J138").Value = 5356.518
$ws.Range("K138").Value = 2777.1177
$ws.Range("L138").Value = 16069.554
$ws.Range("M138").Value = 2362.8823
$ws.Range("N138").Value = -26349.554

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3410.9268
$ws.Range("I32").Value = 3320.4595
$ws.Range("K32").Value = 3320.4595
$ws.Range("M32").Value = -3033.4595
$ws.Range("H61").Value = 4984.846
$ws.Range("I61").Value = 4566.3335
$ws.Range("K61").Value = 4566.3335
$ws.Range("M61").Value = -4354.3335
$ws.Range("H74").Value = 2499.6428
$ws.Range("I74").Value = 2599.04
$ws.Range("K74").Value = 2599.04
$ws.Range("M74").Value = -1725.04
$ws.Range("H77").Value = 2499.6428
$ws.Range("I77").Value = 2599.04
$ws.Range("K77").Value = 12995.2
$ws.Range("M77").Value = -8627.2
$ws.Range("H97").Value = 804
$ws.Range("I97").Value = 875
$ws.Range("J97").Value = 733
$ws.Range("K97").Value = 875
$ws.Range("L97").Value = 733
$ws.Range("M97").Value = -379
$ws.Range("N97").Value = -1725
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents() | Out-Null
$ws.Range("H132").Value = 2604.4
$ws.Range("I132").Value = 2159.0386
$ws.Range("K132").Value = 6477.1158
$ws.Range("M132").Value = -3947.1158
$ws.Range("H136").Value = 4984.846
$ws.Range("I136").Value = 4566.3335
$ws.Range("K136").Value = 13699.0005
$ws.Range("M136").Value = -11149.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1833.5
$ws.Range("I20").Value = 1976.3636
$ws.Range("J20").Value = 1519.2
$ws.Range("K20").Value = 1976.3636
$ws.Range("L20").Value = 1519.2
$ws.Range("M20").Value = -1729.3636
$ws.Range("N20").Value = -2013.2
$ws.Range("H105").Value = 2375.6365
$ws.Range("I105").Value = 2147.75
$ws.Range("K105").Value = 2147.75
$ws.Range("M105").Value = -400.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1719.76
$ws.Range("I58").Value = 1495
$ws.Range("K58").Value = 1495
$ws.Range("M58").Value = -1292
$ws.Range("H99").Value = 3191.6667
$ws.Range("I99").Value = 3175
$ws.Range("K99").Value = 3175
$ws.Range("M99").Value = -1677
$ws.Range("H126").Value = 3191.6667
$ws.Range("I126").Value = 3175
$ws.Range("K126").Value = 9525
$ws.Range("M126").Value = -7055
$ws.Range("H134").Value = 3137.6155
$ws.Range("I134").Value = 3065.75
$ws.Range("K134").Value = 9197.25
$ws.Range("M134").Value = -6662.25
$ws.Range("H136").Value = 1719.76
$ws.Range("I136").Value = 1495
$ws.Range("K136").Value = 4485
$ws.Range("M136").Value = -1935

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12.782609
$ws.Range("I2").Value = 9.6470585
$ws.Range("J2").Value = 21.666666
$ws.Range("K2").Value = 57.882351
$ws.Range("L2").Value = 129.999996
$ws.Range("M2").Value = 55.117649
$ws.Range("N2").Value = -355.999996
$ws.Range("H4").Value = 8790974
$ws.Range("J4").Value = 17000040
$ws.Range("L4").Value = 51000120
$ws.Range("N4").Value = -51000344
$ws.Range("H113").Value = 1886.1818
$ws.Range("I113").Value = 2112.25
$ws.Range("J113").Value = 1757
$ws.Range("K113").Value = 6336.75
$ws.Range("L113").Value = 5271
$ws.Range("M113").Value = -4166.75
$ws.Range("N113").Value = -9611

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1295.375
$ws.Range("I31").Value = 1295.375
$ws.Range("K31").Value = 1295.375
$ws.Range("M31").Value = -1003.375
$ws.Range("H36").Value = 6750
$ws.Range("I36").Value = 5000
$ws.Range("J36").Value = 8500
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 8500
$ws.Range("M36").Value = -4515
$ws.Range("N36").Value = -9470
$ws.Range("H37").Value = 1295.375
$ws.Range("I37").Value = 1295.375
$ws.Range("K37").Value = 1295.375
$ws.Range("M37").Value = -1018.375
$ws.Range("H55").Value = 30939.777
$ws.Range("J55").Value = 34904.668
$ws.Range("L55").Value = 34904.668
$ws.Range("N55").Value = -35558.668
$ws.Range("H70").Value = 30314478
$ws.Range("I70").Value = 66675850
$ws.Range("K70").Value = 66675850
$ws.Range("M70").Value = -66675580
$ws.Range("H73").Value = 30314478
$ws.Range("I73").Value = 66675850
$ws.Range("K73").Value = 66675850
$ws.Range("M73").Value = -66674914
$ws.Range("H97").Value = 849.9375
$ws.Range("I97").Value = 668.1
$ws.Range("J97").Value = 1153
$ws.Range("K97").Value = 668.1
$ws.Range("L97").Value = 1153
$ws.Range("M97").Value = -172.1
$ws.Range("N97").Value = -2145
$ws.Range("H102").Value = 1918.5428
$ws.Range("I102").Value = 1972.9688
$ws.Range("K102").Value = 1972.9688
$ws.Range("M102").Value = -350.9688000000001
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344
$ws.Range("H126").Value = 2597.818
$ws.Range("I126").Value = 2268
$ws.Range("J126").Value = 3175
$ws.Range("K126").Value = 6804
$ws.Range("L126").Value = 9525
$ws.Range("M126").Value = -4334
$ws.Range("N126").Value = -14465
$ws.Range("H132").Value = 3017.8096
$ws.Range("I132").Value = 2626.389
$ws.Range("K132").Value = 7879.167
$ws.Range("M132").Value = -5349.167

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents() | Out-Null
$ws.Range("H7").Value = 5399.1665
$ws.Range("I7").Value = 5299.5
$ws.Range("J7").Value = 5897.5
$ws.Range("K7").Value = 5299.5
$ws.Range("L7").Value = 5897.5
$ws.Range("M7").Value = -5187.5
$ws.Range("N7").Value = -6121.5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents() | Out-Null
$ws.Range("H30").Value = 989.8
$ws.Range("I30").Value = 1199.6666
$ws.Range("J30").Value = 675
$ws.Range("K30").Value = 1199.6666
$ws.Range("L30").Value = 675
$ws.Range("M30").Value = -1091.6666
$ws.Range("N30").Value = -891
$ws.Range("H31").Value = 479
$ws.Range("I31").Value = 200
$ws.Range("J31").Value = 572
$ws.Range("K31").Value = 200
$ws.Range("L31").Value = 572
$ws.Range("M31").Value = 48
$ws.Range("N31").Value = -1068
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents() | Out-Null
$ws.Range("H82").Value = 4107.6
$ws.Range("I82").Value = 2796.5715
$ws.Range("J82").Value = 7166.6665
$ws.Range("K82").Value = 2796.5715
$ws.Range("L82").Value = 7166.6665
$ws.Range("M82").Value = -2435.5715
$ws.Range("N82").Value = -7888.6665
$ws.Range("H85").Value = 4107.6
$ws.Range("I85").Value = 2796.5715
$ws.Range("J85").Value = 7166.6665
$ws.Range("K85").Value = 2796.5715
$ws.Range("L85").Value = 7166.6665
$ws.Range("M85").Value = -1548.5715
$ws.Range("N85").Value = -9662.6665
$ws.Range("H126").Value = 5399.1665
$ws.Range("I126").Value = 5299.5
$ws.Range("J126").Value = 5897.5
$ws.Range("K126").Value = 15898.5
$ws.Range("L126").Value = 17692.5
$ws.Range("M126").Value = -13428.5
$ws.Range("N126").Value = -22632.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4761.1816
$ws.Range("J122").Value = 6726.769
$ws.Range("L122").Value = 20180.307
$ws.Range("N122").Value = -25080.307
$ws.Range("H126").Value = 4667.5884
$ws.Range("J126").Value = 3974
$ws.Range("L126").Value = 11922
$ws.Range("N126").Value = -16862
$ws.Range("H132").Value = 5033.3335
$ws.Range("I132").Value = 4535
$ws.Range("K132").Value = 13605
$ws.Range("M132").Value = -11075
